# Apply a cyclic rotation of data among rows 10, 11, 12 on the active sheet:
#   new row10 = old row11
#   new row11 = old row12
#   new row12 = old row10
# Only the columns that actually differ between these rows are touched:
# A, B, D, E, F, G, H, Q, R, Z, AB

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R", "Z", "AB")

# Capture the original values for rows 10, 11, 12 before making any changes.
$orig10 = @{}
$orig11 = @{}
$orig12 = @{}

foreach ($col in $cols) {
    $orig10[$col] = $ws.Range("$col" + "10").Value2
    $orig11[$col] = $ws.Range("$col" + "11").Value2
    $orig12[$col] = $ws.Range("$col" + "12").Value2
}

foreach ($col in $cols) {
    $ws.Range("$col" + "10").Value = $orig11[$col]
    $ws.Range("$col" + "11").Value = $orig12[$col]
    $ws.Range("$col" + "12").Value = $orig10[$col]
}
